# Update the "want-to-go" counts (column F) on both the "展览" and
# "全部类型" sheets. Both sheets share duplicate rows for these events,
# so both need to be kept in sync.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 674
    8  = 3208
    9  = 4207
    10 = 107
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
